$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grade updates (ДЗ_2 column "D", and a couple of "C" cells) -----------
# Row 7  (Белолипецкий Никита)
$ws.Range("D7").Value = 2

# Row 10 (Горинова Полина)
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 2

# Row 14 (Желтобрюх Максим)
$ws.Range("D14").Value = 2

# Row 15 (Жоркин Игорь)
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 2

# Row 24 (Мельников Сергей)
$ws.Range("D24").Value = 2

# Row 29 (Сафутин Артём)
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 2

# Row 30
$ws.Range("D30").Value = 2

# Recalculate so the shared SUM formulas in column J pick up the new totals
$excel.CalculateFull()

# --- View / scroll position update ----------------------------------------
# Move the frozen pane's visible top-left corner and the active selection,
# keeping the existing freeze split (2 columns / 3 rows) untouched.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 22
$ws.Range("D32").Select()
